$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook holds an "Estado de Cuenta" (account statement) with a table of
# overdue contribution periods for workers of a company. This edit replaces the
# previous multi-worker overdue-period listing with a fresh listing (part 1 of
# a new statement) for a single worker (VERONICA MEDRANO GOMEZ) covering six
# periods (2503-2508).

# Row 32 carries the special "closing" border formatting that marks the last
# row of the table. Row 21 will become the new last row, so bring that
# closing format onto it before the now-unneeded rows are removed.
$ws.Range("B32:J32").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)  # xlPasteFormats

# Remove the old extra worker rows (22:32) - this shifts the trailing
# signature rows (old 37:38) up to become rows 26:27, and keeps the blank
# spacer rows (old 33:36) as the new blank rows 22:25.
$ws.Rows("22:32").Delete()

# Update the summary figures.
$ws.Range("E11").Value = 341640
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 6

# Rewrite the worker detail rows (16:21) for the single worker across the six
# overdue periods.
$periods = @("2503", "2504", "2505", "2506", "2507", "2508")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "1002202673"
    $ws.Cells.Item($r, 4).Value = "VERONICA MEDRANO GOMEZ"
    $ws.Cells.Item($r, 5).Value = $periods[$i]
    $ws.Cells.Item($r, 6).Value = 56940
    $ws.Cells.Item($r, 7).Value = 1423500
}

# Column D ("Nombre Trabajador") is best-fit to its longest content; refresh
# it now that the longest name in the sheet has changed.
$ws.Columns.Item(4).AutoFit()
